# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets
# to reflect the latest scrape, per the gh-pages data refresh commit.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (sheet1) ---
$ws1.Range("F3").Value = 16447
$ws1.Range("F5").Value = 738
$ws1.Range("F6").Value = 15617
$ws1.Range("F7").Value = 71
$ws1.Range("F8").Value = 9285
$ws1.Range("F9").Value = 498
$ws1.Range("F12").Value = 129
$ws1.Range("F13").Value = 227
$ws1.Range("F15").Value = 224
$ws1.Range("F17").Value = 96
$ws1.Range("F18").Value = 623
$ws1.Range("F24").Value = 24
$ws1.Range("F25").Value = 34
$ws1.Range("F26").Value = 537
$ws1.Range("F30").Value = 85
$ws1.Range("F32").Value = 67
$ws1.Range("F33").Value = 271
$ws1.Range("F34").Value = 374
$ws1.Range("F35").Value = 481
$ws1.Range("F37").Value = 5716
$ws1.Range("F38").Value = 5252

# --- Sheet "全部类型" (sheet4) ---
$ws4.Range("F3").Value = 16447
$ws4.Range("F5").Value = 738
$ws4.Range("F6").Value = 15617
$ws4.Range("F7").Value = 71
$ws4.Range("F8").Value = 9285
$ws4.Range("F9").Value = 498
$ws4.Range("F12").Value = 129
$ws4.Range("F13").Value = 227
$ws4.Range("F15").Value = 224
$ws4.Range("F17").Value = 96
$ws4.Range("F18").Value = 623
$ws4.Range("F24").Value = 24
$ws4.Range("F25").Value = 34
$ws4.Range("F26").Value = 537
$ws4.Range("F32").Value = 85
$ws4.Range("F34").Value = 67
$ws4.Range("F35").Value = 271
$ws4.Range("F36").Value = 374
$ws4.Range("F37").Value = 481
$ws4.Range("F39").Value = 5716
$ws4.Range("F41").Value = 5252
